$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 1
$ws.Range("F2").Value = 0.5
$ws.Range("I2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.4
$ws.Range("S2").Value = 1

# Row 3
$ws.Range("C3").Value = 0.75
$ws.Range("F3").Value = 0.75
$ws.Range("I3").Value = 0.5
$ws.Range("M3").Value = 0.5
$ws.Range("S3").Value = 0.5

# Row 4
$ws.Range("C4").Value = 0.8571428571428571
$ws.Range("F4").Value = 0.6
$ws.Range("I4").Value = 0.5714285714285715
$ws.Range("M4").Value = 0.4444444444444445
$ws.Range("S4").Value = 0.6666666666666666

# Row 5
$ws.Range("C5").Value = 0.7894736842105263
$ws.Range("F5").Value = 0.6818181818181818
$ws.Range("I5").Value = 0.5263157894736842
$ws.Range("M5").Value = 0.4761904761904762
$ws.Range("S5").Value = 0.5555555555555556

# Row 6
$ws.Range("C6").Value = 0.7679133988345253
$ws.Range("F6").Value = 0.9639241360211297
$ws.Range("I6").Value = 0.357467732314226
$ws.Range("M6").Value = 0.2250114532569562
$ws.Range("S6").Value = 0.9095582591723214
